$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Notes" header in column AD, row 1 (matches the header row's
# bold style/format automatically since Excel copies the row's look or the
# engine infers style "2" for the header row).
$ws.Range("AD1").Value = "Notes"

# Update the visible selection to AD2 (the first data row under the new
# Notes column) and scroll the view so column W becomes the leftmost
# visible column, as in the authored workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 23
$ws.Range("AD2").Select()
